# Update cryptos list:
#  - Rows 2-33: refresh Price (D) and Volume(1h) (E) figures only.
#  - Rows 34-51: a new coin ("Frax") is inserted at row 34, shifting the
#    previous rows 34-50 down to 35-51 (and the former last row,
#    "Decentraland", drops off the bottom of the 51-row table).
#
# Price values are textual (not real numbers - e.g. "26.588.32",
# "1.744.05", "0.4820" with a significant trailing zero, "1.000", etc.),
# so each Price cell is forced to Text format before the value is
# assigned, to prevent Excel from auto-converting/normalizing the string
# into a floating point number. The cell style is then reset back to
# "Normal" so we don't leave a stray number-format override behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Rows 2-33: update Price (D) and Volume(1h) (E) only ---

$priceVolUpdates = @(
    @{ Row = 2;  D = "26.588.32";      E = "  +3.97%  " },
    @{ Row = 3;  D = "1.744.05";       E = "  +4.45%  " },
    @{ Row = 4;  D = "0.9995";         E = "  +0.04%  " },
    @{ Row = 5;  D = "246.62";         E = "  +3.76%  " },
    @{ Row = 6;  D = $null;            E = "  +0.02%  " },
    @{ Row = 7;  D = "0.4820";         E = "  +0.94%  " },
    @{ Row = 8;  D = "0.2702";         E = "  +3.22%  " },
    @{ Row = 9;  D = "0.06270";        E = "  +1.57%  " },
    @{ Row = 10; D = "1.744.28";       E = "  +4.45%  " },
    @{ Row = 11; D = "0.07138";        E = "  +2.15%  " },
    @{ Row = 12; D = "15.86";          E = "  +7.13%  " },
    @{ Row = 13; D = "0.6246";         E = "  +6.20%  " },
    @{ Row = 14; D = "4.523";          E = "  +3.42%  " },
    @{ Row = 15; D = $null;            E = "  +2.91%  " },
    @{ Row = 17; D = "26.581.19";      E = "  +3.96%  " },
    @{ Row = 18; D = "1.000";          E = "  +0.03%  " },
    @{ Row = 19; D = "0.000006897";    E = "  +2.20%  " },
    @{ Row = 20; D = $null;            E = "  +2.60%  " },
    @{ Row = 21; D = "1.967.43";       E = "  +4.30%  " },
    @{ Row = 22; D = "4.643";          E = "  +4.53%  " },
    @{ Row = 23; D = "8.853";          E = "  +0.73%  " },
    @{ Row = 24; D = "5.378";          E = "  +2.28%  " },
    @{ Row = 25; D = "136.15";         E = "  -0.30%  " },
    @{ Row = 26; D = "15.41";          E = "  +2.26%  " },
    @{ Row = 27; D = "1.819";          E = "  +5.79%  " },
    @{ Row = 28; D = $null;            E = "  +3.50%  " },
    @{ Row = 29; D = "106.89";         E = "  +2.08%  " },
    @{ Row = 30; D = "4.007";          E = "  +0.63%  " },
    @{ Row = 31; D = "3.757";          E = "  +3.55%  " },
    @{ Row = 32; D = "0.07890";        E = "  +0.26%  " },
    @{ Row = 33; D = "0.04631";        E = "  +7.75%  " }
)

foreach ($u in $priceVolUpdates) {
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# --- Rows 34-51: full record shift (B, C, D, E) ---
# A new coin (Frax) is inserted at row 34, pushing the existing rows 34-50
# down to 35-51, and the previous last row (Decentraland, row 51) drops off.

$records = @(
    @{ Row = 34; B = "Frax";             C = "https://coinranking.com/coin/KfWtaeV1W+frax-frax";                     D = "0.9998";      E = "  +0.06%  " },
    @{ Row = 35; B = "HuobiToken";       C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht";              D = "2.619";       E = "  +0.02%  " },
    @{ Row = 36; B = "ImmutableX";       C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";                 D = "0.6447";      E = "  +6.52%  " },
    @{ Row = 37; B = "ARBITRUM";         C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb";                   D = "0.9995";      E = "  +4.75%  " },
    @{ Row = 38; B = "TrustWalletToken"; C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";           D = "0.9383";      E = "  +1.79%  " },
    @{ Row = 39; B = "Quant";            C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt";                  D = "113.52";      E = "  +16.37%  " },
    @{ Row = 40; B = "RenderToken";      C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";           D = "2.002";       E = "  +7.87%  " },
    @{ Row = 41; B = "MXToken";          C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx";                  D = "2.419";       E = "  -6.33%  " },
    @{ Row = 42; B = "PaxDollar";        C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp";                 D = "1.004";       E = "  +0.45%  " },
    @{ Row = 43; B = "FraxShare";        C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs";                  D = "5.757";       E = "  +18.01%  " },
    @{ Row = 44; B = "VeChain";          C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";                D = "0.01510";     E = "  +2.41%  " },
    @{ Row = 45; B = "TheSandbox";       C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand";                D = "0.3924";      E = "  +4.39%  " },
    @{ Row = 46; B = "Algorand";         C = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo";              D = "0.1219";      E = "  +8.94%  " },
    @{ Row = 47; B = "Aptos";            C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt";                     D = "6.769";       E = "  +8.86%  " },
    @{ Row = 48; B = "Cronos";           C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro";                  D = "0.05335";     E = "  +1.30%  " },
    @{ Row = 49; B = "EnergySwap";       C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";                 D = "7.931";       E = "  +6.75%  " },
    @{ Row = 50; B = "Elrond";           C = "https://coinranking.com/coin/omwkOTglq+elrond-egld";                    D = "30.77";       E = "  +2.84%  " },
    @{ Row = 51; B = "NEARProtocol";     C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near";              D = "1.269";       E = "  +5.47%  " }
)

foreach ($r in $records) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
